$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from E1 ("panel") into F1 so the new header cell reuses the
# same bold/border/center style index as the other header cells.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill F2:F112 with the time_taken values captured for each panel gene row.
$values = @(
    "2021-10-05 13:39:32.610699"
    "2021-10-05 13:39:32.610713"
    "2021-10-05 13:39:32.610717"
    "2021-10-05 13:39:32.610720"
    "2021-10-05 13:39:32.610723"
    "2021-10-05 13:39:32.610727"
    "2021-10-05 13:39:32.610729"
    "2021-10-05 13:39:32.610732"
    "2021-10-05 13:39:32.610736"
    "2021-10-05 13:39:32.610739"
    "2021-10-05 13:39:32.610742"
    "2021-10-05 13:39:32.610745"
    "2021-10-05 13:39:32.610747"
    "2021-10-05 13:39:32.610750"
    "2021-10-05 13:39:32.610753"
    "2021-10-05 13:39:32.610756"
    "2021-10-05 13:39:32.610759"
    "2021-10-05 13:39:32.610763"
    "2021-10-05 13:39:32.610766"
    "2021-10-05 13:39:32.610768"
    "2021-10-05 13:39:32.610771"
    "2021-10-05 13:39:32.610774"
    "2021-10-05 13:39:32.610777"
    "2021-10-05 13:39:32.610780"
    "2021-10-05 13:39:32.610783"
    "2021-10-05 13:39:32.610786"
    "2021-10-05 13:39:32.610789"
    "2021-10-05 13:39:32.610792"
    "2021-10-05 13:39:32.610795"
    "2021-10-05 13:39:32.610798"
    "2021-10-05 13:39:32.610801"
    "2021-10-05 13:39:32.610804"
    "2021-10-05 13:39:32.610807"
    "2021-10-05 13:39:32.610810"
    "2021-10-05 13:39:32.610813"
    "2021-10-05 13:39:32.610816"
    "2021-10-05 13:39:32.610819"
    "2021-10-05 13:39:32.610822"
    "2021-10-05 13:39:32.610825"
    "2021-10-05 13:39:32.610827"
    "2021-10-05 13:39:32.610831"
    "2021-10-05 13:39:32.610834"
    "2021-10-05 13:39:32.610837"
    "2021-10-05 13:39:32.610840"
    "2021-10-05 13:39:32.610843"
    "2021-10-05 13:39:32.610846"
    "2021-10-05 13:39:32.610849"
    "2021-10-05 13:39:32.610852"
    "2021-10-05 13:39:32.610855"
    "2021-10-05 13:39:32.610858"
    "2021-10-05 13:39:32.610861"
    "2021-10-05 13:39:32.610864"
    "2021-10-05 13:39:32.610867"
    "2021-10-05 13:39:32.610870"
    "2021-10-05 13:39:32.610873"
    "2021-10-05 13:39:32.610876"
    "2021-10-05 13:39:32.610883"
    "2021-10-05 13:39:32.610886"
    "2021-10-05 13:39:32.610889"
    "2021-10-05 13:39:32.610892"
    "2021-10-05 13:39:32.610895"
    "2021-10-05 13:39:32.610912"
    "2021-10-05 13:39:32.610914"
    "2021-10-05 13:39:32.610917"
    "2021-10-05 13:39:32.610920"
    "2021-10-05 13:39:32.610923"
    "2021-10-05 13:39:32.610926"
    "2021-10-05 13:39:32.610928"
    "2021-10-05 13:39:32.610931"
    "2021-10-05 13:39:32.610933"
    "2021-10-05 13:39:32.610935"
    "2021-10-05 13:39:32.610938"
    "2021-10-05 13:39:32.610940"
    "2021-10-05 13:39:32.610943"
    "2021-10-05 13:39:32.610945"
    "2021-10-05 13:39:32.610948"
    "2021-10-05 13:39:32.610952"
    "2021-10-05 13:39:32.610955"
    "2021-10-05 13:39:32.610957"
    "2021-10-05 13:39:32.610960"
    "2021-10-05 13:39:32.610962"
    "2021-10-05 13:39:32.610965"
    "2021-10-05 13:39:32.610967"
    "2021-10-05 13:39:32.610969"
    "2021-10-05 13:39:32.610972"
    "2021-10-05 13:39:32.610974"
    "2021-10-05 13:39:32.610977"
    "2021-10-05 13:39:32.610979"
    "2021-10-05 13:39:32.610982"
    "2021-10-05 13:39:32.610984"
    "2021-10-05 13:39:32.610986"
    "2021-10-05 13:39:32.610989"
    "2021-10-05 13:39:32.610992"
    "2021-10-05 13:39:32.610995"
    "2021-10-05 13:39:32.610998"
    "2021-10-05 13:39:32.611000"
    "2021-10-05 13:39:32.611002"
    "2021-10-05 13:39:32.611005"
    "2021-10-05 13:39:32.611007"
    "2021-10-05 13:39:32.611010"
    "2021-10-05 13:39:32.611012"
    "2021-10-05 13:39:32.611015"
    "2021-10-05 13:39:32.611017"
    "2021-10-05 13:39:32.611019"
    "2021-10-05 13:39:32.611022"
    "2021-10-05 13:39:32.611024"
    "2021-10-05 13:39:32.611027"
    "2021-10-05 13:39:32.611029"
    "2021-10-05 13:39:32.611033"
    "2021-10-05 13:39:32.611037"
    "2021-10-05 13:39:32.611039"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $values[$i]
}

Write-Output "done"